# Auto update stock data
# For each company's "latest" data row, roll the as-of date forward from
# 2025/11/20 to 2025/11/21 and refresh the EBITDA figure (column B) where it
# changed. Values are kept as plain text (matching how the sheet already
# stores them) rather than letting Excel auto-convert them into a date
# serial number / floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = $originalFormat
}

# Row number => new EBITDA (column B) value. $null means the EBITDA value is
# unchanged for that row, only the date needs to move forward.
$rowUpdates = [ordered]@{
    2  = "4.36"
    8  = "7.50"
    14 = $null
    20 = $null
    26 = "9.58"
    32 = "24.80"
    38 = $null
    44 = "9.73"
    50 = "10.91"
    56 = "29.64"
    62 = "10.54"
    68 = "11.48"
    74 = "14.72"
}

foreach ($r in $rowUpdates.Keys) {
    Set-TextValue $ws.Cells.Item($r, 1) "2025/11/21"

    $newEbitda = $rowUpdates[$r]
    if ($null -ne $newEbitda) {
        Set-TextValue $ws.Cells.Item($r, 2) $newEbitda
    }
}
